$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parvo-ncbi-refseqs-side-data")

# New data widened columns A and B; set explicit widths to match
$ws.Columns.Item(1).ColumnWidth = 13.83
$ws.Columns.Item(2).ColumnWidth = 16.5

# Add new row 26 for HHV-6 Rep (U6)
$ws.Range("A26").Value = "X59532"
$ws.Range("B26").Value = "HHV6-Rep"
$ws.Range("C26").Value = "HHV6-Rep"
$ws.Range("D26").Value = "Parvovirinae"
$ws.Range("E26").Value = "Betherpesparvovirus"
$ws.Range("F26").Value = "not-set"
$ws.Range("G26").Value = "not-set"
$ws.Range("H26").Value = "not-set"
$ws.Range("I26").Value = "not-set"
$ws.Range("J26").Value = "not-set"
$ws.Range("K26").Value = "not-set"
$ws.Range("L26").Value = "not-set"

# Copy style from row 25 so new row matches formatting of preceding data rows
$ws.Range("A25:L25").Copy()
$ws.Range("A26:L26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update row 23 (NC_001539 / CPV / Carnivore protoparvovirus) - set Proto1 clade assignments
$ws.Range("F23").Value = "Proto1"
$ws.Range("H23").Value = "Proto1"
$ws.Range("J23").Value = "Proto1"
$ws.Range("K23").Value = "Proto1"

$ws.Range("C23").Select()
